# PETCARE-v4 deck fix: "Fix a part of a doc. Need more fix."
#
# The only substantive content change between the before/after OOXML is on
# the "4. Use case" slide (slide index 10 in the slide show order): the
# screenshot picture ("Picture 5") is nudged upward so it sits closer to the
# title/divider line (its vertical offset moves from 2289823 EMU to
# 2004073 EMU; horizontal offset and size are unchanged).
#
# Everything else visible in the raw XML diff (stripped creationId extLst
# blocks, dropped cached field text, re-numbered r:id relationship ids,
# endParaRPr bookkeeping, theme identity metadata, etc.) is inert
# round-trip noise with no visible effect, so it is intentionally not
# reproduced here.

$p = $ppt.ActivePresentation

# "4. Use case" is the 10th slide in the slide show order.
$s = $p.Slides.Item(10)

$pic = $s.Shapes.Item("Picture 5")

# Convert the target EMU offset to points (1 pt = 12700 EMU), which is the
# unit PowerPoint's COM object model uses for Shape.Top/Left.
$pic.Top = 2004073 / 12700
